$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Step 1: extend the table with 40 new rows (82-121) for group B2-D.
# First, copy the cell formatting (fills/fonts/number formats) from the
# last two existing data rows (80 = even-row style, 81 = odd-row style)
# down across the new rows so the alternating banding pattern continues
# and no new cell styles are introduced.
# ------------------------------------------------------------------
$ws.Range("A80:G81").Copy()
$ws.Range("A82:G121").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# Step 2: fill in the data for the new B2-D rows.
# Columns A, B, C (Year, Group, Subject), F (Start Time) and G (Duration)
# are plain text/numbers and can be written directly.
# Columns D (Session number) and E (Date) look like numbers/dates, so a
# direct .Value assignment would make Excel silently reinterpret them as
# a real number / date serial. To keep them as literal text (matching the
# source data) without touching cell formatting, we stage each value as a
# text formula in a scratch pair of cells, then copy only the *computed
# values* (PasteSpecial xlPasteValues) into the target cells.
# ------------------------------------------------------------------
$data = @(
    @(82, "Year 5", "B2-D", "endocrinology", "1", "20/12/2025", "09:00:00", 360),
    @(83, "Year 5", "B2-D", "endocrinology", "2", "21/12/2025", "09:00:00", 360),
    @(84, "Year 5", "B2-D", "endocrinology", "3", "22/12/2025", "09:00:00", 360),
    @(85, "Year 5", "B2-D", "endocrinology", "4", "23/12/2025", "09:00:00", 360),
    @(86, "Year 5", "B2-D", "endocrinology", "5", "24/12/2025", "09:00:00", 360),
    @(87, "Year 5", "B2-D", "endocrinology", "6", "27/12/2025", "09:00:00", 360),
    @(88, "Year 5", "B2-D", "endocrinology", "7", "28/12/2025", "09:00:00", 360),
    @(89, "Year 5", "B2-D", "endocrinology", "8", "29/12/2025", "09:00:00", 360),
    @(90, "Year 5", "B2-D", "endocrinology", "9", "30/12/2025", "09:00:00", 360),
    @(91, "Year 5", "B2-D", "endocrinology", "10", "31/12/2025", "09:00:00", 360),
    @(92, "Year 5", "B2-D", "gastroenterology", "1", "03/01/2026", "09:00:00", 360),
    @(93, "Year 5", "B2-D", "gastroenterology", "2", "04/01/2026", "09:00:00", 360),
    @(94, "Year 5", "B2-D", "gastroenterology", "3", "05/01/2026", "09:00:00", 360),
    @(95, "Year 5", "B2-D", "gastroenterology", "4", "06/01/2026", "09:00:00", 360),
    @(96, "Year 5", "B2-D", "gastroenterology", "5", "07/01/2026", "09:00:00", 360),
    @(97, "Year 5", "B2-D", "gastroenterology", "6", "10/01/2026", "09:00:00", 360),
    @(98, "Year 5", "B2-D", "gastroenterology", "7", "11/01/2026", "09:00:00", 360),
    @(99, "Year 5", "B2-D", "gastroenterology", "8", "12/01/2026", "09:00:00", 360),
    @(100, "Year 5", "B2-D", "gastroenterology", "9", "13/01/2026", "09:00:00", 360),
    @(101, "Year 5", "B2-D", "gastroenterology", "10", "14/01/2026", "09:00:00", 360),
    @(102, "Year 5", "B2-D", "neurology", "1", "06/12/2025", "09:00:00", 360),
    @(103, "Year 5", "B2-D", "neurology", "2", "07/12/2025", "09:00:00", 360),
    @(104, "Year 5", "B2-D", "neurology", "3", "08/12/2025", "09:00:00", 360),
    @(105, "Year 5", "B2-D", "neurology", "4", "09/12/2025", "09:00:00", 360),
    @(106, "Year 5", "B2-D", "neurology", "5", "13/12/2025", "09:00:00", 360),
    @(107, "Year 5", "B2-D", "neurology", "6", "14/12/2025", "09:00:00", 360),
    @(108, "Year 5", "B2-D", "neurology", "7", "15/12/2025", "09:00:00", 360),
    @(109, "Year 5", "B2-D", "neurology", "8", "16/12/2025", "09:00:00", 360),
    @(110, "Year 5", "B2-D", "physical medicine", "1", "10/12/2025", "09:00:00", 360),
    @(111, "Year 5", "B2-D", "physical medicine", "2", "17/12/2025", "09:00:00", 360),
    @(112, "Year 5", "B2-D", "rheumatology", "1", "17/01/2026", "09:00:00", 360),
    @(113, "Year 5", "B2-D", "rheumatology", "2", "18/01/2026", "09:00:00", 360),
    @(114, "Year 5", "B2-D", "rheumatology", "3", "19/01/2026", "09:00:00", 360),
    @(115, "Year 5", "B2-D", "rheumatology", "4", "20/01/2026", "09:00:00", 360),
    @(116, "Year 5", "B2-D", "rheumatology", "5", "21/01/2026", "09:00:00", 360),
    @(117, "Year 5", "B2-D", "rheumatology", "6", "07/02/2026", "09:00:00", 360),
    @(118, "Year 5", "B2-D", "rheumatology", "7", "08/02/2026", "09:00:00", 360),
    @(119, "Year 5", "B2-D", "rheumatology", "8", "09/02/2026", "09:00:00", 360),
    @(120, "Year 5", "B2-D", "rheumatology", "9", "10/02/2026", "09:00:00", 360),
    @(121, "Year 5", "B2-D", "rheumatology", "10", "11/02/2026", "09:00:00", 360)
)

foreach ($row in $data) {
    $r       = $row[0]
    $year    = $row[1]
    $group   = $row[2]
    $subject = $row[3]
    $session = $row[4]
    $date    = $row[5]
    $time    = $row[6]
    $dur     = $row[7]

    $ws.Cells.Item($r, 1).Value = $year
    $ws.Cells.Item($r, 2).Value = $group
    $ws.Cells.Item($r, 3).Value = $subject
    $ws.Cells.Item($r, 6).Value = $time
    $ws.Cells.Item($r, 7).Value = $dur

    $ws.Range("Z1").Formula = '="' + $session + '"'
    $ws.Range("AA1").Formula = '="' + $date + '"'
    $ws.Range("Z1:AA1").Copy()
    $ws.Cells.Item($r, 4).Resize(1, 2).PasteSpecial(-4163)
    $excel.CutCopyMode = 0
}

# clean up the scratch cells so they don't linger in the saved sheet
$ws.Range("Z1:AA1").Clear()

Write-Host "Added rows 82-121 (B2-D group)"
